$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040841954743475
$ws.Range("D2").Value = 1.047248882650956
$ws.Range("E2").Value = 1.048779703637209
$ws.Range("F2").Value = 1.058631582101418
$ws.Range("I2").Value = 1.040657333725013
$ws.Range("J2").Value = 1.045926379842882
$ws.Range("K2").Value = 1.05001204938615
$ws.Range("L2").Value = 1.051538596173364
$ws.Range("M2").Value = 1.061363287630382

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041789827902005
$ws.Range("D3").Value = 1.047963768525999
$ws.Range("E3").Value = 1.049596995860929
$ws.Range("F3").Value = 1.059473834301504
$ws.Range("I3").Value = 1.040844806873969
$ws.Range("J3").Value = 1.046520044721692
$ws.Range("K3").Value = 1.050538992793409
$ws.Range("L3").Value = 1.052167988974688
$ws.Range("M3").Value = 1.062019539487099

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042403625761958
$ws.Range("D4").Value = 1.048426429402844
$ws.Range("E4").Value = 1.050126540764277
$ws.Range("F4").Value = 1.060019336896839
$ws.Range("I4").Value = 1.040964521540481
$ws.Range("J4").Value = 1.046904020301928
$ws.Range("K4").Value = 1.050879369848693
$ws.Range("L4").Value = 1.052575298091038
$ws.Range("M4").Value = 1.062444039145143

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04266177549961
$ws.Range("D5").Value = 1.048620950089612
$ws.Range("E5").Value = 1.050349328026034
$ws.Range("F5").Value = 1.060248786453723
$ws.Range("I5").Value = 1.041014467468612
$ws.Range("J5").Value = 1.047065402951568
$ws.Range("K5").Value = 1.051022321784938
$ws.Range("L5").Value = 1.052746541571358
$ws.Range("M5").Value = 1.062622464135098

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042705126342132
$ws.Range("D6").Value = 1.048653611982235
$ws.Range("E6").Value = 1.050386744711946
$ws.Range("F6").Value = 1.060287319049401
$ws.Range("I6").Value = 1.041022831180278
$ws.Range("J6").Value = 1.047092497390433
$ws.Range("K6").Value = 1.051046315648409
$ws.Range("L6").Value = 1.052775294698919
$ws.Range("M6").Value = 1.062652420403944

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042407074744051
$ws.Range("D7").Value = 1.048429028528803
$ws.Range("E7").Value = 1.050129517005922
$ws.Range("F7").Value = 1.060022402342518
$ws.Range("I7").Value = 1.040965190422823
$ws.Range("J7").Value = 1.046906176866051
$ws.Range("K7").Value = 1.050881280541202
$ws.Range("L7").Value = 1.05257758621418
$ws.Range("M7").Value = 1.062446423406503

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041162196872237
$ws.Range("D8").Value = 1.04749046378489
$ws.Range("E8").Value = 1.049055765066056
$ws.Range("F8").Value = 1.058916118793409
$ws.Range("I8").Value = 1.040721020543168
$ws.Range("J8").Value = 1.046127045224978
$ws.Range("K8").Value = 1.050190253940472
$ws.Range("L8").Value = 1.051751291357824
$ws.Range("M8").Value = 1.061585099026019

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038972140361892
$ws.Range("D9").Value = 1.045837288638574
$ws.Range("E9").Value = 1.047169121826968
$ws.Range("F9").Value = 1.056970676060698
$ws.Range("I9").Value = 1.040278594578323
$ws.Range("J9").Value = 1.0447528987659
$ws.Range("K9").Value = 1.048968106353969
$ws.Range("L9").Value = 1.050295689857731
$ws.Range("M9").Value = 1.060066328824306

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037514574166634
$ws.Range("D10").Value = 1.044735731347859
$ws.Range("E10").Value = 1.045915110091231
$ws.Range("F10").Value = 1.055676478678548
$ws.Range("I10").Value = 1.03997550721785
$ws.Range("J10").Value = 1.043836048199366
$ws.Range("K10").Value = 1.04815040581109
$ws.Range("L10").Value = 1.049325653209698
$ws.Range("M10").Value = 1.059053217797845

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036884031659322
$ws.Range("D11").Value = 1.044258896679017
$ws.Range("E11").Value = 1.045373017398822
$ws.Range("F11").Value = 1.055116752421304
$ws.Range("I11").Value = 1.039842346216686
$ws.Range("J11").Value = 1.043438876721925
$ws.Range("K11").Value = 1.047795650753816
$ws.Range("L11").Value = 1.048905717903949
$ws.Range("M11").Value = 1.058614404270923

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036649910269895
$ws.Range("D12").Value = 1.044081802637775
$ws.Range("E12").Value = 1.045171797138243
$ws.Range("F12").Value = 1.054908947362323
$ws.Range("I12").Value = 1.039792596201591
$ws.Range("J12").Value = 1.043291325199106
$ws.Range("K12").Value = 1.047663777107247
$ws.Range("L12").Value = 1.048749750972154
$ws.Range("M12").Value = 1.058451391189444

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036700126005447
$ws.Range("D13").Value = 1.0441197888192
$ws.Range("E13").Value = 1.045214953339195
$ws.Range("F13").Value = 1.054953517603145
$ws.Range("I13").Value = 1.039803280774266
$ws.Range("J13").Value = 1.043322976586635
$ws.Range("K13").Value = 1.047692069041521
$ws.Range("L13").Value = 1.048783205677498
$ws.Range("M13").Value = 1.058486358851555

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036864677266718
$ws.Range("D14").Value = 1.044244257540436
$ws.Range("E14").Value = 1.045356381663881
$ws.Range("F14").Value = 1.055099573096479
$ws.Range("I14").Value = 1.039838239737861
$ws.Range("J14").Value = 1.04342668055768
$ws.Range("K14").Value = 1.047784752107625
$ws.Range("L14").Value = 1.048892825300891
$ws.Range("M14").Value = 1.05860092992253

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036966074709573
$ws.Range("D15").Value = 1.044320949998159
$ws.Range("E15").Value = 1.045443538534229
$ws.Range("F15").Value = 1.055189576295141
$ws.Range("I15").Value = 1.039859740955118
$ws.Range("J15").Value = 1.043490572801649
$ws.Range("K15").Value = 1.04784184374511
$ws.Range("L15").Value = 1.048960367689238
$ws.Range("M15").Value = 1.05867151856782

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037556433784521
$ws.Range("D16").Value = 1.044767380518971
$ws.Range("E16").Value = 1.045951106142974
$ws.Range("F16").Value = 1.055713640111979
$ws.Range("I16").Value = 1.039984304236276
$ws.Range("J16").Value = 1.043862403689211
$ws.Range("K16").Value = 1.048173935394194
$ws.Range("L16").Value = 1.049353525075809
$ws.Range("M16").Value = 1.059082337783127

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037926910004966
$ws.Range("D17").Value = 1.045047454992076
$ws.Range("E17").Value = 1.046269732613674
$ws.Range("F17").Value = 1.056042551920834
$ws.Range("I17").Value = 1.040061925362006
$ws.Range("J17").Value = 1.044095598847425
$ws.Range("K17").Value = 1.048382065036125
$ws.Range("L17").Value = 1.049600169104447
$ws.Range("M17").Value = 1.059339999997136

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038143059688851
$ws.Range("D18").Value = 1.045210831766843
$ws.Range("E18").Value = 1.046455669037769
$ws.Range("F18").Value = 1.056234465130057
$ws.Range("I18").Value = 1.040107014942993
$ws.Range("J18").Value = 1.044231601248224
$ws.Range("K18").Value = 1.04850339733468
$ws.Range("L18").Value = 1.049744041670713
$ws.Range("M18").Value = 1.059490277399844

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038216770811076
$ws.Range("D19").Value = 1.045266541365896
$ws.Range("E19").Value = 1.046519083263848
$ws.Range("F19").Value = 1.056299913470041
$ws.Range("I19").Value = 1.040122357838289
$ws.Range("J19").Value = 1.044277971740098
$ws.Range("K19").Value = 1.048544757249544
$ws.Range("L19").Value = 1.049793100035746
$ws.Range("M19").Value = 1.059541515932726

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037887155491493
$ws.Range("D20").Value = 1.045017404185098
$ws.Range("E20").Value = 1.046235537999811
$ws.Range("F20").Value = 1.056007256111972
$ws.Range("I20").Value = 1.040053616532782
$ws.Range("J20").Value = 1.044070580901474
$ws.Range("K20").Value = 1.048359741534415
$ws.Range("L20").Value = 1.049573705566322
$ws.Range("M20").Value = 1.059312356574716

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036816218524036
$ws.Range("D21").Value = 1.044207603941107
$ws.Range("E21").Value = 1.045314730750654
$ws.Range("F21").Value = 1.055056560560341
$ws.Range("I21").Value = 1.039827953141543
$ws.Range("J21").Value = 1.043396142986706
$ws.Range("K21").Value = 1.047757462056392
$ws.Range("L21").Value = 1.048860544610152
$ws.Range("M21").Value = 1.05856719208973

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036143400175186
$ws.Range("D22").Value = 1.04369858759383
$ws.Range("E22").Value = 1.044736576650717
$ws.Range("F22").Value = 1.054459412577566
$ws.Range("I22").Value = 1.039684403152175
$ws.Range("J22").Value = 1.04297195673072
$ws.Range("K22").Value = 1.047378196490866
$ws.Range("L22").Value = 1.048412243647202
$ws.Range("M22").Value = 1.058098572946278

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036500024063011
$ws.Range("D23").Value = 1.043968413280169
$ws.Range("E23").Value = 1.045042991323698
$ws.Range("F23").Value = 1.054775915378547
$ws.Range("I23").Value = 1.03976065944711
$ws.Range("J23").Value = 1.043196838843471
$ws.Range("K23").Value = 1.047579307751326
$ws.Range("L23").Value = 1.048649887393461
$ws.Range("M23").Value = 1.058347006308754

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03790511866833
$ws.Range("D24").Value = 1.045030982807922
$ws.Range("E24").Value = 1.046250988805281
$ws.Range("F24").Value = 1.056023204569523
$ws.Range("I24").Value = 1.040057371508307
$ws.Range("J24").Value = 1.0440818854851
$ws.Range("K24").Value = 1.048369828768587
$ws.Range("L24").Value = 1.049585663271538
$ws.Range("M24").Value = 1.059324847486489

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039537891707631
$ws.Range("D25").Value = 1.046264582120025
$ws.Range("E25").Value = 1.047656209838127
$ws.Range("F25").Value = 1.060459078438381
$ws.Range("I25").Value = 1.040394409659634
$ws.Range("J25").Value = 1.045108286045295
$ws.Range("K25").Value = 1.049284583325093
$ws.Range("L25").Value = 1.050671938740634
$ws.Range("M25").Value = 1.060459078438381
